# [ADDITIONAL SCRAPING] add a "Player Info" sheet and an "ODI Batting Extra"
# sheet, and replace the MATCH_CARD_LINK url columns on the existing sheets
# with a plain MATCH_CODE column (just the numeric code out of the url).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New first sheet: "Player Info"
# ---------------------------------------------------------------------------
$battingSheetForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetForInsert, $null)
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $playerInfo.Cells.Item(1, $c).Value = $piHeaders[$c - 1]
}
$piHeaderRange = $playerInfo.Range("A1:D1")
$piHeaderRange.Font.Bold = $true
$piHeaderRange.Borders.LineStyle = 1
$piHeaderRange.HorizontalAlignment = -4108
$piHeaderRange.VerticalAlignment = -4160

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "3766"
$playerInfo.Cells.Item(2, 2).Value = "Daniel Trevor Christian"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Fast Medium"

# ---------------------------------------------------------------------------
# 2) "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE, url -> bare code
#    NOTE: re-fetch the sheet handle by name now that a sheet was inserted
#    in front of it - worksheet handles obtained before an Add() that
#    shifts their slot do not track the sheet they originally pointed to.
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

for ($r = 2; $r -le 21; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $url = $cell.Value()
    $code = $url.Substring($url.LastIndexOf("=") + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $code
}

# these two rows had a stray blank INNING_NUMBER cell that should be removed
# entirely (no cell at all) rather than left as an empty string
$battingSheet.Range("B19").ClearContents()
$battingSheet.Range("B21").ClearContents()

# ---------------------------------------------------------------------------
# 3) "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE, url -> bare code
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

for ($r = 2; $r -le 20; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $url = $cell.Value()
    $code = $url.Substring($url.LastIndexOf("=") + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $code
}

# ---------------------------------------------------------------------------
# 4) New last sheet: "ODI Batting Extra"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $exHeaders.Length; $c++) {
    $extra.Cells.Item(1, $c).Value = $exHeaders[$c - 1]
}
$exHeaderRange = $extra.Range("A1:F1")
$exHeaderRange.Font.Bold = $true
$exHeaderRange.Borders.LineStyle = 1
$exHeaderRange.HorizontalAlignment = -4108
$exHeaderRange.VerticalAlignment = -4160

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$exData = @(
    @("3367", 7, "0", "0", "7.87%", "NO"),
    @("3371", 7, "4", "1", "14.29%", "NO"),
    @("3373", $null, $null, $null, $null, "NO"),
    @("3377", 7, "1", "0", "3.80%", "NO"),
    @("3380", 7, "5", "0", "10.42%", "NO"),
    @("3384", 7, "0", "0", "2.14%", "NO"),
    @("3386", 7, "2", "0", "9.52%", "NO"),
    @("3389", 7, "0", "0", "1.31%", "NO"),
    @("3391", 4, "2", "0", "3.12%", "NO"),
    @("3392", 7, "0", "0", "1.48%", "NO"),
    @("3393", 7, "2", "0", "8.23%", "NO"),
    @("3398", $null, $null, $null, $null, "NO"),
    @("3400", 7, "0", "0", "3.90%", "NO"),
    @("3402", 7, "0", "1", "5.45%", "NO"),
    @("3439", 8, "0", "0", "1.51%", "NO"),
    @("3441", 8, "2", "0", "7.26%", "NO"),
    @("3443", 8, "0", "0", "0.80%", "NO"),
    @("3602", 7, $null, $null, $null, "NO"),
    @("3606", $null, $null, $null, $null, "NO"),
    @("4486", $null, $null, $null, $null, "NO")
)

for ($i = 0; $i -lt $exData.Length; $i++) {
    $row = $i + 2
    $rowData = $exData[$i]

    $cellA = $extra.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $rowData[0]

    $cellB = $extra.Cells.Item($row, 2)
    if ($rowData[1] -ne $null) {
        $cellB.Value = $rowData[1]
    }

    $cellC = $extra.Cells.Item($row, 3)
    if ($rowData[2] -ne $null) {
        $cellC.NumberFormat = "@"
        $cellC.Value = $rowData[2]
    }

    $cellD = $extra.Cells.Item($row, 4)
    if ($rowData[3] -ne $null) {
        $cellD.NumberFormat = "@"
        $cellD.Value = $rowData[3]
    }

    $cellE = $extra.Cells.Item($row, 5)
    if ($rowData[4] -ne $null) {
        $cellE.NumberFormat = "@"
        $cellE.Value = $rowData[4]
    }

    $cellF = $extra.Cells.Item($row, 6)
    $cellF.Value = $rowData[5]
}
